# Tunnel source for cellular = DONE
#
# - "Main Info" sheet, B7 (Design): BASE -> SMART
# - "Main Info" sheet, B25 (4G+Cellular, backup link): FALSE -> TRUE
# - Selection on "Main Info" sheet moves from D15 to D10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# Design drop-down value for the tunnel: BASE -> SMART
$ws.Range("B7").Value = "SMART"

# 4G+Cellular flag used as backup (tunnel source for cellular): False -> True
$ws.Range("B25").Value = $true

# Move the active cell selection
$ws.Activate()
$ws.Range("D10").Select()
